$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table_PredatorWeightAtAge")

# Update the data values for rows 3-27 (columns B-I) with the revised figures
$ws.Cells.Item(3, 2).Value = 695
$ws.Cells.Item(3, 3).Value = 695
$ws.Cells.Item(3, 4).Value = 28.04426553
$ws.Cells.Item(3, 5).Value = 27.85136412
$ws.Cells.Item(3, 6).Value = 34.10662519
$ws.Cells.Item(3, 7).Value = 26.56700203
$ws.Cells.Item(3, 8).Value = 95.55793084
$ws.Cells.Item(3, 9).Value = 96.8855249
$ws.Cells.Item(4, 2).Value = 949
$ws.Cells.Item(4, 3).Value = 949
$ws.Cells.Item(4, 4).Value = 40.37419797
$ws.Cells.Item(4, 5).Value = 37.64431759
$ws.Cells.Item(4, 6).Value = 50.48576957
$ws.Cells.Item(4, 7).Value = 46.83028492
$ws.Cells.Item(4, 8).Value = 133.4273667
$ws.Cells.Item(4, 9).Value = 125.8213997
$ws.Cells.Item(5, 2).Value = 1208
$ws.Cells.Item(5, 3).Value = 1208
$ws.Cells.Item(5, 4).Value = 50.04673738
$ws.Cells.Item(5, 5).Value = 45.30089613
$ws.Cells.Item(5, 6).Value = 73.05964723
$ws.Cells.Item(5, 7).Value = 62.15104386
$ws.Cells.Item(5, 8).Value = 182.0900859
$ws.Cells.Item(5, 9).Value = 152.0571917
$ws.Cells.Item(6, 2).Value = 1455
$ws.Cells.Item(6, 3).Value = 1455
$ws.Cells.Item(6, 4).Value = 57.63461501
$ws.Cells.Item(6, 5).Value = 51.28715912
$ws.Cells.Item(6, 6).Value = 102.5661152
$ws.Cells.Item(6, 7).Value = 72.30517212
$ws.Cells.Item(6, 8).Value = 241.4814331
$ws.Cells.Item(6, 9).Value = 174.9160102
$ws.Cells.Item(7, 2).Value = 1682
$ws.Cells.Item(7, 3).Value = 1682
$ws.Cells.Item(7, 4).Value = 63.58712494
$ws.Cells.Item(7, 5).Value = 55.96749308
$ws.Cells.Item(7, 6).Value = 138.5739924
$ws.Cells.Item(7, 7).Value = 78.61185388
$ws.Cells.Item(7, 8).Value = 309.5719141
$ws.Cells.Item(7, 9).Value = 194.2857141
$ws.Cells.Item(8, 2).Value = 1881
$ws.Cells.Item(8, 3).Value = 1881
$ws.Cells.Item(8, 4).Value = 68.25672773
$ws.Cells.Item(8, 5).Value = 59.62679205
$ws.Cells.Item(8, 6).Value = 179.0074464
$ws.Cells.Item(8, 7).Value = 82.3966125
$ws.Cells.Item(8, 8).Value = 382.2571594
$ws.Cells.Item(8, 9).Value = 210.3735324
$ws.Cells.Item(9, 2).Value = 2051
$ws.Cells.Item(9, 3).Value = 2051
$ws.Cells.Item(9, 4).Value = 71.9199203
$ws.Cells.Item(9, 5).Value = 62.48779921
$ws.Cells.Item(9, 6).Value = 220.3794251
$ws.Cells.Item(9, 7).Value = 84.62548892
$ws.Cells.Item(9, 8).Value = 454.1669864
$ws.Cells.Item(9, 9).Value = 223.5406004
$ws.Cells.Item(10, 2).Value = 2194
$ws.Cells.Item(10, 3).Value = 2194
$ws.Cells.Item(10, 4).Value = 74.79360795
$ws.Cells.Item(10, 5).Value = 64.72466535
$ws.Cells.Item(10, 6).Value = 258.8740244
$ws.Cells.Item(10, 7).Value = 85.92430606
$ws.Cells.Item(10, 8).Value = 520.1529661
$ws.Cells.Item(10, 9).Value = 234.1998294
$ws.Cells.Item(11, 2).Value = 2311
$ws.Cells.Item(11, 3).Value = 2311
$ws.Cells.Item(11, 4).Value = 77.04794771
$ws.Cells.Item(11, 5).Value = 66.47354951
$ws.Cells.Item(11, 6).Value = 291.6509561
$ws.Cells.Item(11, 7).Value = 86.67664462
$ws.Cells.Item(11, 8).Value = 576.6515475
$ws.Cells.Item(11, 9).Value = 242.7579796
$ws.Cells.Item(12, 2).Value = 2406
$ws.Cells.Item(12, 3).Value = 2406
$ws.Cells.Item(12, 4).Value = 78.81642376
$ws.Cells.Item(12, 5).Value = 67.84090707
$ws.Cells.Item(12, 6).Value = 317.5165414
$ws.Cells.Item(12, 7).Value = 87.11095374
$ws.Cells.Item(12, 8).Value = 622.228393
$ws.Cells.Item(12, 9).Value = 249.5862734
$ws.Cells.Item(13, 2).Value = 2482
$ws.Cells.Item(13, 3).Value = 2482
$ws.Cells.Item(13, 4).Value = 80.20375125
$ws.Cells.Item(13, 5).Value = 68.9099697
$ws.Cells.Item(13, 6).Value = 336.7305302
$ws.Cells.Item(13, 7).Value = 87.36118313
$ws.Cells.Item(13, 8).Value = 657.2596521
$ws.Cells.Item(13, 9).Value = 255.0083227
$ws.Cells.Item(14, 2).Value = 2726
$ws.Cells.Item(14, 3).Value = 2589
$ws.Cells.Item(14, 4).Value = 81.29207675
$ws.Cells.Item(14, 5).Value = 69.74581171
$ws.Cells.Item(14, 6).Value = 350.3718541
$ws.Cells.Item(14, 7).Value = 87.50519318
$ws.Cells.Item(14, 8).Value = 683.1978866
$ws.Cells.Item(14, 9).Value = 259.2979137
$ws.Cells.Item(15, 2).Value = 2970
$ws.Cells.Item(15, 3).Value = 2696
$ws.Cells.Item(15, 4).Value = 82.14584227
$ws.Cells.Item(15, 5).Value = 70.39931119
$ws.Cells.Item(15, 6).Value = 359.7485363
$ws.Cells.Item(15, 7).Value = 87.58801962
$ws.Cells.Item(15, 8).Value = 701.8768305
$ws.Cells.Item(15, 9).Value = 262.6819525
$ws.Cells.Item(16, 2).Value = 3214
$ws.Cells.Item(16, 3).Value = 2803
$ws.Cells.Item(16, 4).Value = 82.81560106
$ws.Cells.Item(16, 5).Value = 70.91024695
$ws.Cells.Item(16, 6).Value = 366.0514041
$ws.Cells.Item(16, 7).Value = 87.6356392
$ws.Cells.Item(16, 8).Value = 715.0605081
$ws.Cells.Item(16, 9).Value = 265.3457501
$ws.Cells.Item(17, 2).Value = 3458
$ws.Cells.Item(17, 3).Value = 2910
$ws.Cells.Item(17, 4).Value = 83.34101094
$ws.Cells.Item(17, 5).Value = 71.30971987
$ws.Cells.Item(17, 6).Value = 370.2247147
$ws.Cells.Item(17, 7).Value = 87.66301143
$ws.Cells.Item(17, 8).Value = 724.2341683
$ws.Cells.Item(17, 9).Value = 267.4390336
$ws.Cells.Item(18, 2).Value = 3702
$ws.Cells.Item(18, 3).Value = 3017
$ws.Cells.Item(18, 4).Value = 83.75318253
$ws.Cells.Item(18, 5).Value = 71.62204606
$ws.Cells.Item(18, 6).Value = 372.9604719
$ws.Cells.Item(18, 7).Value = 87.67874337
$ws.Cells.Item(18, 8).Value = 730.5545085
$ws.Cells.Item(18, 9).Value = 269.0818195
$ws.Cells.Item(19, 2).Value = 3946
$ws.Cells.Item(19, 3).Value = 3124
$ws.Cells.Item(19, 4).Value = 84.07652135
$ws.Cells.Item(19, 5).Value = 71.86623695
$ws.Cells.Item(19, 6).Value = 374.7421154
$ws.Cells.Item(19, 7).Value = 87.68778453
$ws.Cells.Item(19, 8).Value = 734.8793026
$ws.Cells.Item(19, 9).Value = 270.3697362
$ws.Cells.Item(20, 2).Value = 4190
$ws.Cells.Item(20, 3).Value = 3231
$ws.Cells.Item(20, 4).Value = 84.33017299
$ws.Cells.Item(20, 5).Value = 72.05715654
$ws.Cells.Item(20, 6).Value = 375.8974386
$ws.Cells.Item(20, 7).Value = 87.69298028
$ws.Cells.Item(20, 8).Value = 737.8247686
$ws.Cells.Item(20, 9).Value = 271.3786344
$ws.Cells.Item(21, 2).Value = 4434
$ws.Cells.Item(21, 3).Value = 3338
$ws.Cells.Item(21, 4).Value = 84.52915668
$ws.Cells.Item(21, 5).Value = 72.2064262
$ws.Cells.Item(21, 6).Value = 376.6445387
$ws.Cells.Item(21, 7).Value = 87.6959661
$ws.Cells.Item(21, 8).Value = 739.8244216
$ws.Cells.Item(21, 9).Value = 272.1684694
$ws.Cells.Item(22, 2).Value = 4434
$ws.Cells.Item(22, 3).Value = 3338
$ws.Cells.Item(22, 4).Value = 84.68525466
$ws.Cells.Item(22, 5).Value = 72.32313205
$ws.Cells.Item(22, 6).Value = 377.1267894
$ws.Cells.Item(22, 7).Value = 87.69768192
$ws.Cells.Item(22, 8).Value = 741.1790261
$ws.Cells.Item(22, 9).Value = 272.7865067
$ws.Cells.Item(23, 2).Value = 4434
$ws.Cells.Item(23, 3).Value = 3338
$ws.Cells.Item(23, 4).Value = 84.80770982
$ws.Cells.Item(23, 5).Value = 72.41437801
$ws.Cells.Item(23, 6).Value = 377.4377195
$ws.Cells.Item(23, 7).Value = 87.69866792
$ws.Cells.Item(23, 8).Value = 742.095313
$ws.Cells.Item(23, 9).Value = 273.2699312
$ws.Cells.Item(24, 2).Value = 4434
$ws.Cells.Item(24, 3).Value = 3338
$ws.Cells.Item(24, 4).Value = 84.90377299
$ws.Cells.Item(24, 5).Value = 72.48571827
$ws.Cells.Item(24, 6).Value = 377.6380408
$ws.Cells.Item(24, 7).Value = 87.69923452
$ws.Cells.Item(24, 8).Value = 742.7144949
$ws.Cells.Item(24, 9).Value = 273.6479508
$ws.Cells.Item(25, 2).Value = 4434
$ws.Cells.Item(25, 3).Value = 3338
$ws.Cells.Item(25, 4).Value = 84.97913227
$ws.Cells.Item(25, 5).Value = 72.54149535
$ws.Cells.Item(25, 6).Value = 377.7670386
$ws.Cells.Item(25, 7).Value = 87.69956012
$ws.Cells.Item(25, 8).Value = 743.1326264
$ws.Cells.Item(25, 9).Value = 273.9434796
$ws.Cells.Item(26, 2).Value = 4434
$ws.Cells.Item(26, 3).Value = 3338
$ws.Cells.Item(26, 4).Value = 85.03824983
$ws.Cells.Item(26, 5).Value = 72.58510442
$ws.Cells.Item(26, 6).Value = 377.8500814
$ws.Cells.Item(26, 7).Value = 87.69974722
$ws.Cells.Item(26, 8).Value = 743.4148609
$ws.Cells.Item(26, 9).Value = 274.1744772
$ws.Cells.Item(27, 2).Value = 4434
$ws.Cells.Item(27, 3).Value = 3338
$ws.Cells.Item(27, 4).Value = 85.08462616
$ws.Cells.Item(27, 5).Value = 72.61919998
$ws.Cells.Item(27, 6).Value = 377.9035298
$ws.Cells.Item(27, 7).Value = 87.69985474
$ws.Cells.Item(27, 8).Value = 743.605308
$ws.Cells.Item(27, 9).Value = 274.3550091

# Apply number formatting (0.00) and right alignment to the D3:I27 block
$ws.Range("D3:I27").NumberFormat = "0.00"
$ws.Range("D3:I27").HorizontalAlignment = -4152

# Make this sheet the active/selected tab, with D4 as the active cell
$ws.Activate() | Out-Null
$ws.Range("D4").Select() | Out-Null

Write-Host "edit complete"
